# Insert a new data row above current row 107, shifting existing rows 107-161
# down to 108-162 (so the sheet grows from 161 to 162 rows incl. header).
# The new row represents a "Papa" record for "Región del Maule", 1a (cosecha),
# dated 2022-01-11 (serial 44572).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 107 (and everything below it) down by one row, carrying
# formatting (e.g. the date number format on column D) along with it.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted, now-empty row 107 with its data.
$ws.Range("A107").Value = 11
$ws.Range("B107").Value = "Vega Monumental Concepción"
$ws.Range("C107").Value = "Bíobío"
$ws.Range("D107").Value = 44572
$ws.Range("E107").Value = 8
$ws.Range("F107").Value = 100114001
$ws.Range("G107").Value = "Papa"
$ws.Range("H107").Value = "Asterix"
$ws.Range("I107").Value = "1a (cosecha)"
$ws.Range("J107").Value = 320
$ws.Range("K107").Value = 12000
$ws.Range("L107").Value = 13000
$ws.Range("M107").Value = 12531
$ws.Range("N107").Value = "$/saco 25 kilos"
$ws.Range("O107").Value = "Región del Maule"
$ws.Range("P107").Value = 501
$ws.Range("Q107").Value = 25
$ws.Range("R107").Value = "Hortaliza"
